$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the full target range is treated as text so numeric-looking strings are not converted to numbers
$ws.Range("A1:K13").NumberFormat = "@"

$data = New-Object 'object[,]' 13,11

$data[0,0] = "venue"
$data[0,1] = "date"
$data[0,2] = "result"
$data[0,3] = "ownTeam"
$data[0,4] = "oppTeam"
$data[0,5] = "batsman"
$data[0,6] = "totalRuns"
$data[0,7] = "totalBalls"
$data[0,8] = "total4s"
$data[0,9] = "total6s"
$data[0,10] = "sr"

$data[1,0] = " Dubai (DSC)"
$data[1,1] = " October 17 2020"
$data[1,2] = "RCB won by 7 wickets (with 2 balls remaining)"
$data[1,3] = "Royal Challengers Bangalore"
$data[1,4] = "Rajasthan Royals"
$data[1,5] = "Aaron Finch "
$data[1,6] = "14"
$data[1,7] = "11"
$data[1,8] = "0"
$data[1,9] = "2"
$data[1,10] = "127.27"

$data[2,0] = " Abu Dhabi"
$data[2,1] = " October 03 2020"
$data[2,2] = "RCB won by 8 wickets (with 5 balls remaining)"
$data[2,3] = "Royal Challengers Bangalore"
$data[2,4] = "Rajasthan Royals"
$data[2,5] = "Aaron Finch "
$data[2,6] = "8"
$data[2,7] = "7"
$data[2,8] = "2"
$data[2,9] = "0"
$data[2,10] = "114.28"

$data[3,0] = " Dubai (DSC)"
$data[3,1] = " September 24 2020"
$data[3,2] = "Kings XI won by 97 runs"
$data[3,3] = "Royal Challengers Bangalore"
$data[3,4] = "Kings XI Punjab"
$data[3,5] = "Aaron Finch "
$data[3,6] = "20"
$data[3,7] = "21"
$data[3,8] = "3"
$data[3,9] = "0"
$data[3,10] = "95.23"

$data[4,0] = " Sharjah"
$data[4,1] = " October 15 2020"
$data[4,2] = "Kings XI won by 8 wickets"
$data[4,3] = "Royal Challengers Bangalore"
$data[4,4] = "Kings XI Punjab"
$data[4,5] = "Aaron Finch "
$data[4,6] = "20"
$data[4,7] = "18"
$data[4,8] = "2"
$data[4,9] = "1"
$data[4,10] = "111.11"

$data[5,0] = " Dubai (DSC)"
$data[5,1] = " September 21 2020"
$data[5,2] = "RCB won by 10 runs"
$data[5,3] = "Royal Challengers Bangalore"
$data[5,4] = "Sunrisers Hyderabad"
$data[5,5] = "Aaron Finch "
$data[5,6] = "29"
$data[5,7] = "27"
$data[5,8] = "1"
$data[5,9] = "2"
$data[5,10] = "107.40"

$data[6,0] = " Abu Dhabi"
$data[6,1] = " November 06 2020"
$data[6,2] = "Sunrisers won by 6 wickets (with 2 balls remaining)"
$data[6,3] = "Royal Challengers Bangalore"
$data[6,4] = "Sunrisers Hyderabad"
$data[6,5] = "Aaron Finch "
$data[6,6] = "32"
$data[6,7] = "30"
$data[6,8] = "3"
$data[6,9] = "1"
$data[6,10] = "106.66"

$data[7,0] = " Dubai (DSC)"
$data[7,1] = " October 05 2020"
$data[7,2] = "Capitals won by 59 runs"
$data[7,3] = "Royal Challengers Bangalore"
$data[7,4] = "Delhi Capitals"
$data[7,5] = "Aaron Finch "
$data[7,6] = "13"
$data[7,7] = "14"
$data[7,8] = "1"
$data[7,9] = "0"
$data[7,10] = "92.85"

$data[8,0] = " Abu Dhabi"
$data[8,1] = " October 21 2020"
$data[8,2] = "RCB won by 8 wickets (with 39 balls remaining)"
$data[8,3] = "Royal Challengers Bangalore"
$data[8,4] = "Kolkata Knight Riders"
$data[8,5] = "Aaron Finch "
$data[8,6] = "16"
$data[8,7] = "21"
$data[8,8] = "2"
$data[8,9] = "0"
$data[8,10] = "76.19"

$data[9,0] = " Dubai (DSC)"
$data[9,1] = " September 28 2020"
$data[9,2] = "Match tied (RCB won the one-over eliminator)"
$data[9,3] = "Royal Challengers Bangalore"
$data[9,4] = "Mumbai Indians"
$data[9,5] = "Aaron Finch "
$data[9,6] = "52"
$data[9,7] = "35"
$data[9,8] = "7"
$data[9,9] = "1"
$data[9,10] = "148.57"

$data[10,0] = " Sharjah"
$data[10,1] = " October 12 2020"
$data[10,2] = "RCB won by 82 runs"
$data[10,3] = "Royal Challengers Bangalore"
$data[10,4] = "Kolkata Knight Riders"
$data[10,5] = "Aaron Finch "
$data[10,6] = "47"
$data[10,7] = "37"
$data[10,8] = "4"
$data[10,9] = "1"
$data[10,10] = "127.02"

$data[11,0] = " Dubai (DSC)"
$data[11,1] = " October 25 2020"
$data[11,2] = "Super Kings won by 8 wickets (with 8 balls remaining)"
$data[11,3] = "Royal Challengers Bangalore"
$data[11,4] = "Chennai Super Kings"
$data[11,5] = "Aaron Finch "
$data[11,6] = "15"
$data[11,7] = "11"
$data[11,8] = "3"
$data[11,9] = "0"
$data[11,10] = "136.36"

$data[12,0] = " Dubai (DSC)"
$data[12,1] = " October 10 2020"
$data[12,2] = "RCB won by 37 runs"
$data[12,3] = "Royal Challengers Bangalore"
$data[12,4] = "Chennai Super Kings"
$data[12,5] = "Aaron Finch "
$data[12,6] = "2"
$data[12,7] = "9"
$data[12,8] = "0"
$data[12,9] = "0"
$data[12,10] = "22.22"

$ws.Range("A1:K13").Value = $data

